# Updates crypto price/volume figures (and the Stacks/FirstDigitalUSD row swap)
# per the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-detects numeric-looking text (e.g. "0.999") and coerces it to a
# real number, which would drop the fixed-width formatting the sheet relies on.
# Writing it with a leading apostrophe forces text, then resetting the style back
# to Normal clears the quote-prefix formatting Excel applied while keeping the
# value as text (matching the original inlineStr cells).
function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '56.569.24'
$ws.Cells.Item(2, 5).Value = '  +10.88%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.250.58'
$ws.Cells.Item(3, 5).Value = '  +6.51%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
Set-TextValue 5 4 '398.43'
$ws.Cells.Item(5, 5).Value = '  +3.09%  '

# Row 6
Set-TextValue 6 4 '111.29'
$ws.Cells.Item(6, 5).Value = '  +9.36%  '

# Row 7
Set-TextValue 7 4 '0.560'
$ws.Cells.Item(7, 5).Value = '  +4.78%  '

# Row 8
Set-TextValue 8 4 '0.999'

# Row 9
Set-TextValue 9 4 '0.620'
$ws.Cells.Item(9, 5).Value = '  +7.34%  '

# Row 10
Set-TextValue 10 4 '39.33'
$ws.Cells.Item(10, 5).Value = '  +7.45%  '

# Row 11
Set-TextValue 11 4 '0.0948'
$ws.Cells.Item(11, 5).Value = '  +12.02%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.47%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '3.764.63'
$ws.Cells.Item(13, 5).Value = '  +6.61%  '

# Row 14
Set-TextValue 14 4 '19.20'
$ws.Cells.Item(14, 5).Value = '  +5.20%  '

# Row 15
Set-TextValue 15 4 '8.11'
$ws.Cells.Item(15, 5).Value = '  +5.91%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.246.58'
$ws.Cells.Item(16, 5).Value = '  +6.38%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +6.86%  '

# Row 18
Set-TextValue 18 4 '11.18'
$ws.Cells.Item(18, 5).Value = '  +5.01%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '56.410.30'
$ws.Cells.Item(19, 5).Value = '  +10.47%  '

# Row 20
Set-TextValue 20 4 '3.32'
$ws.Cells.Item(20, 5).Value = '  +3.90%  '

# Row 21
Set-TextValue 21 4 '0.0000104'
$ws.Cells.Item(21, 5).Value = '  +9.67%  '

# Row 22
Set-TextValue 22 4 '13.02'
$ws.Cells.Item(22, 5).Value = '  +6.88%  '

# Row 23
Set-TextValue 23 4 '300.12'
$ws.Cells.Item(23, 5).Value = '  +13.90%  '

# Row 24
Set-TextValue 24 4 '75.51'
$ws.Cells.Item(24, 5).Value = '  +8.78%  '

# Row 25
Set-TextValue 25 4 '3.23'
$ws.Cells.Item(25, 5).Value = '  +3.90%  '

# Row 26
Set-TextValue 26 4 '8.15'
$ws.Cells.Item(26, 5).Value = '  +3.59%  '

# Row 27
Set-TextValue 27 4 '28.32'
$ws.Cells.Item(27, 5).Value = '  +5.30%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +3.93%  '

# Row 29
Set-TextValue 29 4 '7.36'
$ws.Cells.Item(29, 5).Value = '  +2.50%  '

# Row 30
Set-TextValue 30 4 '0.171'
$ws.Cells.Item(30, 5).Value = '  +5.68%  '

# Row 31
Set-TextValue 31 4 '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.12%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +6.80%  '

# Row 33
Set-TextValue 33 4 '11.11'
$ws.Cells.Item(33, 5).Value = '  +7.63%  '

# Row 34
Set-TextValue 34 4 '36.70'
$ws.Cells.Item(34, 5).Value = '  +3.79%  '

# Row 35
Set-TextValue 35 4 '0.0485'
$ws.Cells.Item(35, 5).Value = '  +3.47%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +3.52%  '

# Row 37
Set-TextValue 37 4 '51.58'
$ws.Cells.Item(37, 5).Value = '  +3.15%  '

# Row 38
Set-TextValue 38 4 '3.55'
$ws.Cells.Item(38, 5).Value = '  +6.59%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 39 4 '3.11'
$ws.Cells.Item(39, 5).Value = '  +26.34%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 40 4 '0.999'
$ws.Cells.Item(40, 5).Value = '  -0.07%  '

# Row 41
Set-TextValue 41 4 '17.60'
$ws.Cells.Item(41, 5).Value = '  +7.65%  '

# Row 42
Set-TextValue 42 4 '134.50'
$ws.Cells.Item(42, 5).Value = '  +3.31%  '

# Row 43
Set-TextValue 43 4 '1.93'
$ws.Cells.Item(43, 5).Value = '  +6.21%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +4.98%  '

# Row 45
Set-TextValue 45 4 '3.98'
$ws.Cells.Item(45, 5).Value = '  +6.38%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -2.67%  '

# Row 47
Set-TextValue 47 4 '22.32'
$ws.Cells.Item(47, 5).Value = '  +3.80%  '

# Row 48
Set-TextValue 48 4 '2.20'
$ws.Cells.Item(48, 5).Value = '  +58.64%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.144.37'
$ws.Cells.Item(49, 5).Value = '  +4.78%  '

# Row 50
Set-TextValue 50 4 '2.09'
$ws.Cells.Item(50, 5).Value = '  +0.90%  '

# Row 51
Set-TextValue 51 4 '2.41'
$ws.Cells.Item(51, 5).Value = '  -3.20%  '
